# Fix quantity of potentiometers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nomenclature")

# Update the quantity of the "Potentiomètre 10k" row (row 13) from 1 to 2.
$ws.Range("D13").Value = 2

# Reflect the resulting active cell/selection seen in the saved file.
$ws.Range("D26").Select()

$wb.Save()
